# Scheduled-runner profit refresh: rewrite the currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on each class sheet with refreshed market-board figures.
# Values come straight from the updated data pull; no formulas involved in this sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 384.85
$ws.Range("I2").Value = 413.13333
$ws.Range("K2").Value = 413.13333
$ws.Range("M2").Value = -300.13333
# Row 69
$ws.Range("H69").Value = 5309.7
$ws.Range("I69").Value = 6171
$ws.Range("J69").Value = 3300
$ws.Range("K69").Value = 18513
$ws.Range("L69").Value = 9900
$ws.Range("M69").Value = -17639
$ws.Range("N69").Value = -11648
# Row 72
$ws.Range("H72").Value = 5309.7
$ws.Range("I72").Value = 6171
$ws.Range("J72").Value = 3300
$ws.Range("K72").Value = 55539
$ws.Range("L72").Value = 29700
$ws.Range("M72").Value = -51171
$ws.Range("N72").Value = -38436
# Row 98
$ws.Range("H98").Value = 349.92
$ws.Range("I98").Value = 349.92
$ws.Range("K98").Value = 349.92
$ws.Range("M98").Value = 1148.08
# Row 113
$ws.Range("H113").Value = 1987.5
$ws.Range("I113").Value = 1987.5
$ws.Range("K113").Value = 1987.5
$ws.Range("M113").Value = 1266.5
# Row 122
$ws.Range("H122").Value = 349.92
$ws.Range("I122").Value = 349.92
$ws.Range("K122").Value = 1049.76
$ws.Range("M122").Value = 1400.24
# Row 133
$ws.Range("H133").Value = 98568.42999999999
$ws.Range("J133").Value = 98568.42999999999
$ws.Range("L133").Value = 98568.42999999999
$ws.Range("N133").Value = -108688.43
# Row 134
$ws.Range("H134").Value = 82355.75
$ws.Range("J134").Value = 82355.75
$ws.Range("L134").Value = 82355.75
$ws.Range("N134").Value = -92495.75
# Row 136
$ws.Range("H136").Value = 78605
$ws.Range("J136").Value = 78605
$ws.Range("L136").Value = 78605
$ws.Range("N136").Value = -88805
# Row 137
$ws.Range("H137").Value = 455286.44
$ws.Range("I137").Value = 1369.2084
$ws.Range("J137").Value = 1817038.1
$ws.Range("K137").Value = 4107.6252
$ws.Range("L137").Value = 5451114.300000001
$ws.Range("M137").Value = -1557.6252
$ws.Range("N137").Value = -5456214.300000001
# Row 139
$ws.Range("H139").Value = 69518
$ws.Range("J139").Value = 69518
$ws.Range("L139").Value = 69518
$ws.Range("N139").Value = -79798
# Row 140
$ws.Range("H140").Value = 91557.10000000001
$ws.Range("J140").Value = 91557.10000000001
$ws.Range("L140").Value = 91557.10000000001
$ws.Range("N140").Value = -101917.1

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5725.9834
$ws.Range("I32").Value = 2827.6047
$ws.Range("J32").Value = 13057.177
$ws.Range("K32").Value = 2827.6047
$ws.Range("L32").Value = 13057.177
$ws.Range("M32").Value = -2540.6047
$ws.Range("N32").Value = -13631.177
# Row 43
$ws.Range("H43").Value = 71850
$ws.Range("I43").Value = 55555
$ws.Range("K43").Value = 55555
$ws.Range("M43").Value = -55242
# Row 63
$ws.Range("H63").Value = 4141.222
$ws.Range("I63").Value = 3712
$ws.Range("J63").Value = 4999.6665
$ws.Range("K63").Value = 3712
$ws.Range("L63").Value = 4999.6665
$ws.Range("M63").Value = -3026
$ws.Range("N63").Value = -6371.6665
# Row 66
$ws.Range("H66").Value = 4141.222
$ws.Range("I66").Value = 3712
$ws.Range("J66").Value = 4999.6665
$ws.Range("K66").Value = 18560
$ws.Range("L66").Value = 24998.3325
$ws.Range("M66").Value = -15128
$ws.Range("N66").Value = -31862.3325
# Row 80
$ws.Range("H80").Value = 82703
$ws.Range("J80").Value = 82703
$ws.Range("L80").Value = 82703
$ws.Range("N80").Value = -84699
# Row 83
$ws.Range("H83").Value = 82703
$ws.Range("J83").Value = 82703
$ws.Range("L83").Value = 248109
$ws.Range("N83").Value = -258093

$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 28307.885
$ws.Range("J132").Value = 28307.885
$ws.Range("L132").Value = 28307.885
$ws.Range("N132").Value = -38427.88499999999
# Row 134
$ws.Range("H134").Value = 3714.8386
$ws.Range("I134").Value = 2366.4
$ws.Range("J134").Value = 9333.333000000001
$ws.Range("K134").Value = 7099.200000000001
$ws.Range("L134").Value = 27999.999
$ws.Range("M134").Value = -4564.200000000001
$ws.Range("N134").Value = -33069.999
# Row 135
$ws.Range("H135").Value = 97665.71000000001
$ws.Range("J135").Value = 97665.71000000001
$ws.Range("L135").Value = 97665.71000000001
$ws.Range("N135").Value = -107805.71
# Row 138
$ws.Range("H138").Value = 82370.375
$ws.Range("J138").Value = 82370.375
$ws.Range("L138").Value = 82370.375
$ws.Range("N138").Value = -92650.375
# Row 140
$ws.Range("H140").Value = 75711.28999999999
$ws.Range("J140").Value = 75711.28999999999
$ws.Range("L140").Value = 75711.28999999999
$ws.Range("N140").Value = -86071.28999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 4860.4287
$ws.Range("I58").Value = 5922.5
$ws.Range("J58").Value = 3444.3333
$ws.Range("K58").Value = 5922.5
$ws.Range("L58").Value = 3444.3333
$ws.Range("M58").Value = -5719.5
$ws.Range("N58").Value = -3850.3333
# Row 86
$ws.Range("H86").Value = 2993830.5
$ws.Range("I86").Value = 7156897
$ws.Range("J86").Value = 20211.857
$ws.Range("K86").Value = 7156897
$ws.Range("L86").Value = 20211.857
$ws.Range("M86").Value = -7155774
$ws.Range("N86").Value = -22457.857
# Row 89
$ws.Range("H89").Value = 2993830.5
$ws.Range("I89").Value = 7156897
$ws.Range("J89").Value = 20211.857
$ws.Range("K89").Value = 35784485
$ws.Range("L89").Value = 101059.285
$ws.Range("M89").Value = -35778869
$ws.Range("N89").Value = -112291.285
# Row 99
$ws.Range("H99").Value = 15875613
$ws.Range("I99").Value = 37038716
$ws.Range("J99").Value = 3285
$ws.Range("K99").Value = 37038716
$ws.Range("L99").Value = 3285
$ws.Range("M99").Value = -37037218
$ws.Range("N99").Value = -6281
# Row 122
$ws.Range("H122").Value = 3723.318
$ws.Range("I122").Value = 3438.625
$ws.Range("J122").Value = 3886
$ws.Range("K122").Value = 10315.875
$ws.Range("L122").Value = 11658
$ws.Range("M122").Value = -7865.875
$ws.Range("N122").Value = -16558
# Row 126
$ws.Range("H126").Value = 15875613
$ws.Range("I126").Value = 37038716
$ws.Range("J126").Value = 3285
$ws.Range("K126").Value = 111116148
$ws.Range("L126").Value = 9855
$ws.Range("M126").Value = -111113678
$ws.Range("N126").Value = -14795
# Row 134
$ws.Range("H134").Value = 2553657
$ws.Range("I134").Value = 3403876
$ws.Range("K134").Value = 10211628
$ws.Range("M134").Value = -10209093
# Row 136
$ws.Range("H136").Value = 4860.4287
$ws.Range("I136").Value = 5922.5
$ws.Range("J136").Value = 3444.3333
$ws.Range("K136").Value = 17767.5
$ws.Range("L136").Value = 10332.9999
$ws.Range("M136").Value = -15217.5
$ws.Range("N136").Value = -15432.9999
# Row 138
$ws.Range("H138").Value = 92455.625
$ws.Range("J138").Value = 92455.625
$ws.Range("L138").Value = 92455.625
$ws.Range("N138").Value = -102735.625

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 267.33334
$ws.Range("I13").Value = 184.33333
$ws.Range("K13").Value = 552.99999
$ws.Range("M13").Value = -384.99999
# Row 121
$ws.Range("H121").Value = 668898.6
$ws.Range("I121").Value = 1743.5
$ws.Range("J121").Value = 1431361.6
$ws.Range("K121").Value = 5230.5
$ws.Range("L121").Value = 4294084.800000001
$ws.Range("M121").Value = -3920.5
$ws.Range("N121").Value = -4296704.800000001
# Row 128
$ws.Range("H128").Value = 297498.5
$ws.Range("I128").Value = 297498.5
$ws.Range("K128").Value = 892495.5
$ws.Range("M128").Value = -887515.5

$ws = $wb.Worksheets.Item("GSM")
# Row 109
$ws.Range("H109").Value = 86997
$ws.Range("J109").Value = 86997
$ws.Range("L109").Value = 86997
$ws.Range("N109").Value = -89077
# Row 135
$ws.Range("H135").Value = 99817.60000000001
$ws.Range("J135").Value = 99817.60000000001
$ws.Range("L135").Value = 99817.60000000001
$ws.Range("N135").Value = -109957.6
# Row 140
$ws.Range("H140").Value = 90251.75
$ws.Range("J140").Value = 90251.75
$ws.Range("L140").Value = 90251.75
$ws.Range("N140").Value = -100611.75

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3539603.8
$ws.Range("I40").Value = 88090.914
$ws.Range("K40").Value = 88090.914
$ws.Range("M40").Value = -87954.914
# Row 46
$ws.Range("H46").Value = 8753.571
$ws.Range("I46").Value = 11385.1
$ws.Range("J46").Value = 2174.75
$ws.Range("K46").Value = 11385.1
$ws.Range("L46").Value = 2174.75
$ws.Range("M46").Value = -11197.1
$ws.Range("N46").Value = -2550.75
# Row 100
$ws.Range("H100").Value = 10310.412
$ws.Range("I100").Value = 10377.071
$ws.Range("K100").Value = 10377.071
$ws.Range("M100").Value = -9836.071
# Row 132
$ws.Range("H132").Value = 2833.1667
$ws.Range("I132").Value = 1999.75
$ws.Range("K132").Value = 5999.25
$ws.Range("M132").Value = -3469.25

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 18560
$ws.Range("J45").Value = 18560
$ws.Range("L45").Value = 18560
$ws.Range("N45").Value = -19542
# Row 122
$ws.Range("H122").Value = 3259.7273
$ws.Range("I122").Value = 3345.2778
$ws.Range("K122").Value = 10035.8334
$ws.Range("M122").Value = -7585.8334
# Row 126
$ws.Range("H126").Value = 3501.9285
$ws.Range("I126").Value = 3262.7
$ws.Range("K126").Value = 9788.099999999999
$ws.Range("M126").Value = -7318.099999999999
